$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Junio de 2020 a las 15:29"

# Update country stats (values refreshed + re-sorted rows by Casos totales desc)
# Row 4
$ws.Cells.Item(4, 2).Value = 2162864
$ws.Cells.Item(4, 3).Value = 636
$ws.Cells.Item(4, 5).Value = 1174922
$ws.Cells.Item(4, 7).Value = 7
$ws.Cells.Item(4, 8).Value = 117865

# Row 7
$ws.Cells.Item(7, 2).Value = 333475
$ws.Cells.Item(7, 3).Value = 692
$ws.Cells.Item(7, 5).Value = 153964

# Row 19
$ws.Cells.Item(19, 2).Value = 132048
$ws.Cells.Item(19, 3).Value = 4507
$ws.Cells.Item(19, 4).Value = 87890
$ws.Cells.Item(19, 5).Value = 43147
$ws.Cells.Item(19, 7).Value = 39
$ws.Cells.Item(19, 8).Value = 1011

# Row 23
$ws.Cells.Item(23, 2).Value = 80876
$ws.Cells.Item(23, 3).Value = 1274
$ws.Cells.Item(23, 4).Value = 58681
$ws.Cells.Item(23, 5).Value = 22119
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(23, 8).Value = 76

# Row 38
$ws.Cells.Item(38, 4).Value = 9891
$ws.Cells.Item(38, 5).Value = 20844
$ws.Cells.Item(38, 7).Value = 9
$ws.Cells.Item(38, 8).Value = 842

# Row 47
$ws.Cells.Item(47, 1).Value = "Irak"
$ws.Cells.Item(47, 2).Value = 21315
$ws.Cells.Item(47, 3).Value = 1106
$ws.Cells.Item(47, 4).Value = 9271
$ws.Cells.Item(47, 5).Value = 11392
$ws.Cells.Item(47, 7).Value = 45
$ws.Cells.Item(47, 8).Value = 652

# Row 48
$ws.Cells.Item(48, 1).Value = "Panama"
$ws.Cells.Item(48, 2).Value = 20686
$ws.Cells.Item(48, 4).Value = 13766
$ws.Cells.Item(48, 5).Value = 6483
$ws.Cells.Item(48, 8).Value = 437

# Row 50
$ws.Cells.Item(50, 5).Value = 5301
$ws.Cells.Item(50, 7).Value = 4
$ws.Cells.Item(50, 8).Value = 46

# Row 57
$ws.Cells.Item(57, 2).Value = 12367
$ws.Cells.Item(57, 3).Value = 57
$ws.Cells.Item(57, 5).Value = 601
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(57, 8).Value = 255

# Row 69
$ws.Cells.Item(69, 2).Value = 8639
$ws.Cells.Item(69, 3).Value = 8
$ws.Cells.Item(69, 5).Value = 259

# Row 76
$ws.Cells.Item(76, 2).Value = 5154
$ws.Cells.Item(76, 3).Value = 74
$ws.Cells.Item(76, 5).Value = 1139

# Row 87
$ws.Cells.Item(87, 2).Value = 3727
$ws.Cells.Item(87, 3).Value = 133
$ws.Cells.Item(87, 4).Value = 1286
$ws.Cells.Item(87, 5).Value = 2337
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = 104

# Row 88
$ws.Cells.Item(88, 1).Value = "Etiopia"
$ws.Cells.Item(88, 2).Value = 3521
$ws.Cells.Item(88, 3).Value = 176
$ws.Cells.Item(88, 4).Value = 620
$ws.Cells.Item(88, 5).Value = 2841
$ws.Cells.Item(88, 7).Value = 3
$ws.Cells.Item(88, 8).Value = 60

# Row 89
$ws.Cells.Item(89, 1).Value = "Gabon"
$ws.Cells.Item(89, 2).Value = 3463
$ws.Cells.Item(89, 4).Value = 1024
$ws.Cells.Item(89, 5).Value = 2416
$ws.Cells.Item(89, 8).Value = 23

# Row 93
$ws.Cells.Item(93, 4).Value = 2162
$ws.Cells.Item(93, 5).Value = 713
$ws.Cells.Item(93, 7).Value = 2
$ws.Cells.Item(93, 8).Value = 165

# Row 98
$ws.Cells.Item(98, 1).Value = "Cuba"
$ws.Cells.Item(98, 2).Value = 2262
$ws.Cells.Item(98, 3).Value = 14
$ws.Cells.Item(98, 4).Value = 1965
$ws.Cells.Item(98, 5).Value = 213
$ws.Cells.Item(98, 8).Value = 84

# Row 99
$ws.Cells.Item(99, 1).Value = "Croacia"
$ws.Cells.Item(99, 2).Value = 2254
$ws.Cells.Item(99, 3).Value = 2
$ws.Cells.Item(99, 4).Value = 2140
$ws.Cells.Item(99, 5).Value = 7
$ws.Cells.Item(99, 8).Value = 107

# Row 104
$ws.Cells.Item(104, 1).Value = "Mali"
$ws.Cells.Item(104, 2).Value = 1860
$ws.Cells.Item(104, 3).Value = 51
$ws.Cells.Item(104, 4).Value = 1125
$ws.Cells.Item(104, 5).Value = 631
$ws.Cells.Item(104, 8).Value = 104

# Row 105
$ws.Cells.Item(105, 1).Value = "Islandia"
$ws.Cells.Item(105, 2).Value = 1810
$ws.Cells.Item(105, 4).Value = 1796
$ws.Cells.Item(105, 5).Value = 4
$ws.Cells.Item(105, 8).Value = 10

# Row 122
$ws.Cells.Item(122, 2).Value = 1176
$ws.Cells.Item(122, 3).Value = 7
$ws.Cells.Item(122, 4).Value = 683
$ws.Cells.Item(122, 5).Value = 442

# Row 124
$ws.Cells.Item(124, 1).Value = "Tunez"
$ws.Cells.Item(124, 2).Value = 1110
$ws.Cells.Item(124, 3).Value = 14
$ws.Cells.Item(124, 4).Value = 999
$ws.Cells.Item(124, 5).Value = 62
$ws.Cells.Item(124, 8).Value = 49

# Row 125
$ws.Cells.Item(125, 1).Value = "Letonia"
$ws.Cells.Item(125, 2).Value = 1097
$ws.Cells.Item(125, 4).Value = 845
$ws.Cells.Item(125, 5).Value = 224
$ws.Cells.Item(125, 8).Value = 28

# Row 198
$ws.Cells.Item(198, 1).Value = "Belice"
$ws.Cells.Item(198, 3).Value = 1
$ws.Cells.Item(198, 4).Value = 16
$ws.Cells.Item(198, 5).Value = 3
$ws.Cells.Item(198, 8).Value = 2

# Row 199
$ws.Cells.Item(199, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(199, 2).Value = 21
$ws.Cells.Item(199, 4).Value = 20
$ws.Cells.Item(199, 5).Value = 1
$ws.Cells.Item(199, 8).Value = 0

# Row 206
$ws.Cells.Item(206, 1).Value = "Groenlandia"

# Row 207
$ws.Cells.Item(207, 1).Value = "Islas Malvinas"

# Row 208
$ws.Cells.Item(208, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 8).Value = 1

# Row 209
$ws.Cells.Item(209, 1).Value = "Santa Sede"
$ws.Cells.Item(209, 4).Value = 12
$ws.Cells.Item(209, 8).Value = 0
